# Insert two new daily price rows ("Fruta / hortaliza, semanal") into the
# Limon / Vega Modelo de Temuco sheet. The existing rows 881-921 shift down
# to 883-923 and the two freshly inserted rows (881, 882) get new data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 881:921 down by two rows.
$ws.Rows("881:882").Insert()

# Common (constant) columns for this product block.
$mercadoId = 10
$mercado   = "Vega Modelo de Temuco"
$region    = "La Araucanía"
$codreg    = 9
$tipo      = "Fruta"
$prodId    = 100102
$producto  = "Cítricos"
$catId     = 100102003
$categoria = "Limón"

# --- New row 881 ---
$ws.Cells.Item(881, 1).Value  = $mercadoId
$ws.Cells.Item(881, 2).Value  = $mercado
$ws.Cells.Item(881, 3).Value  = $region
$ws.Cells.Item(881, 4).Value  = 44509
$ws.Cells.Item(881, 5).Value  = $codreg
$ws.Cells.Item(881, 6).Value  = $tipo
$ws.Cells.Item(881, 7).Value  = $prodId
$ws.Cells.Item(881, 8).Value  = $producto
$ws.Cells.Item(881, 9).Value  = $catId
$ws.Cells.Item(881, 10).Value = $categoria
$ws.Cells.Item(881, 11).Value = "Sin especificar"
$ws.Cells.Item(881, 12).Value = "1a amarillo"
$ws.Cells.Item(881, 13).Value = 295
$ws.Cells.Item(881, 14).Value = 8000
$ws.Cells.Item(881, 15).Value = 9000
$ws.Cells.Item(881, 16).Value = 8627
$ws.Cells.Item(881, 17).Value = "$/bandeja 15 kilos"
$ws.Cells.Item(881, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(881, 19).Value = 575
$ws.Cells.Item(881, 20).Value = 15

# --- New row 882 ---
$ws.Cells.Item(882, 1).Value  = $mercadoId
$ws.Cells.Item(882, 2).Value  = $mercado
$ws.Cells.Item(882, 3).Value  = $region
$ws.Cells.Item(882, 4).Value  = 44421
$ws.Cells.Item(882, 5).Value  = $codreg
$ws.Cells.Item(882, 6).Value  = $tipo
$ws.Cells.Item(882, 7).Value  = $prodId
$ws.Cells.Item(882, 8).Value  = $producto
$ws.Cells.Item(882, 9).Value  = $catId
$ws.Cells.Item(882, 10).Value = $categoria
$ws.Cells.Item(882, 11).Value = "Sin especificar"
$ws.Cells.Item(882, 12).Value = "1a amarillo"
$ws.Cells.Item(882, 13).Value = 260
$ws.Cells.Item(882, 14).Value = 8000
$ws.Cells.Item(882, 15).Value = 8000
$ws.Cells.Item(882, 16).Value = 8000
$ws.Cells.Item(882, 17).Value = "$/malla 18 kilos"
$ws.Cells.Item(882, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(882, 19).Value = 444
$ws.Cells.Item(882, 20).Value = 18

# Match the date-number format already used in column D for the rest of
# the sheet (style index 2 -> numFmt 165).
$ws.Range("D881:D882").NumberFormat = $ws.Range("D880").NumberFormat
